# AutenticazioneViaJar.docx — remove "numero di telefono" (phone number)
# field from the registration use case, replacing it with "data di
# nascita" (date of birth), per the commit:
#   "Removed telephone number (modified Domain Model, Authentication SSD
#    and Use Case)"

$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceOne = 1 (not used directly — constants below
# are passed positionally to Find.Execute as in the COM signature):
#   Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#           MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#           Format, ReplaceWith, Replace)

# 1) "Il Sistema richiede l'inserimento di un nome, un cognome ed un
#     numero di telefono." -> "... ed una data di nascita."
$p = $d.Paragraphs.Item(33)
$p.Range.Find.Execute("un numero di telefono", $true, $false, $false, `
    $false, $false, $true, 1, $false, "una data di nascita", 2)

# 2) "L'Utente inserisce il nome, il cognome ed il numero di telefono."
#    -> "... ed il ..." becomes "... e la data di nascita."
$p = $d.Paragraphs.Item(34)
$p.Range.Find.Execute("ed il numero di telefono", $true, $false, $false, `
    $false, $false, $true, 1, $false, "e la data di nascita", 2)

# 3) Paragraph "Il Sistema verifica i dati inseriti dall'Utente." loses
#    its justified alignment (now left/default).
$p = $d.Paragraphs.Item(35)
$p.Format.Alignment = 0

# 4) "L'Utente inserisce un nome che non rispetta il formato corretto:"
#    — text unchanged, runs are coalesced (re-save normalises formatting).
$p = $d.Paragraphs.Item(36)
$p.Range.Find.Execute("L'Utente inserisce un nome che non rispetta il formato corretto:", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "L'Utente inserisce un nome che non rispetta il formato corretto:", 2)

# 5) "L'Utente inserisce un cognome che non rispetta il formato
#    corretto:" — text unchanged, runs are coalesced.
$p = $d.Paragraphs.Item(39)
$p.Range.Find.Execute("L'Utente inserisce un cognome che non rispetta il formato corretto:", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "L'Utente inserisce un cognome che non rispetta il formato corretto:", 2)

# 6) "L'Utente inserisce un numero di telefono che non rispetta il
#    formato corretto:" -> "... una data di nascita ..." and the
#    paragraph loses its justified alignment.
$p = $d.Paragraphs.Item(42)
$p.Range.Find.Execute("un numero di telefono", $true, $false, $false, `
    $false, $false, $true, 1, $false, "una data di nascita", 2)
$p.Format.Alignment = 0

# 7) "Il nome, il cognome ed il numero di telefono devono rispettare il
#    formato definito dal Sistema." -> "... e la data di nascita ..."
$p = $d.Paragraphs.Item(127)
$p.Range.Find.Execute("ed il numero di telefono", $true, $false, $false, `
    $false, $false, $true, 1, $false, "e la data di nascita", 2)
